$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: append a trailing data row (old test entry) ---
$c120a = $ws1.Cells.Item(120, 1)
$c120a.NumberFormat = "@"
$c120a.Value = "00000000"
$c120a.ClearFormats()
$ws1.Cells.Item(120, 4).Value = 2

# --- Add Sheet2 right after Sheet1; it becomes the active sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "serial num"
$ws2.Range("B1").Value = "Current"
$ws2.Range("C1").Value = "Voltage"
$ws2.Range("D1").Value = "CC"
$ws2.Range("E1").Value = "CV"

# Header styling: Calibri, blue, 11pt for most headers
$headerRange = $ws2.Range("A1:E1")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Font.Color = 12611584

# "Current" header uses Consolas 10pt blue, vertically centered
$ws2.Range("B1").Font.Name = "Consolas"
$ws2.Range("B1").Font.Size = 10
$ws2.Range("B1").Font.Color = 12611584
$ws2.Range("B1").VerticalAlignment = -4108

$ws2.Range("A2").Value = 67362781

$ws1.Range("A25").Select()
$ws2.Range("F2").Select()
